$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "t7"
$ws.Range("B8").Value = "w"
$ws.Range("C8").Value = "r"
$ws.Range("D8").Value = "n"

$ws.Range("F5").Select()
